$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Price" (D) / "Volume(1h)" (E) figures from the latest cryptos feed.
# Values are stored as plain scraped text (not native Excel numbers), e.g. a
# price like "26.448.70" uses a thousands-separator dot and is unambiguous
# text, but a price like "1.007" parses as a valid number and Excel would
# normally auto-convert it under General formatting. $priceIsNumericLike
# flags rows where that auto-conversion has to be headed off by briefly
# switching the cell to Text format, writing the literal digits, then
# clearing that temporary format so the cell is left in its original,
# unstyled state, same as every other cell in the sheet.
$updates = @(
    @{ Row = 2; D = "26.448.70"; E = "  -3.42%  " }
    @{ Row = 3; D = "1.804.18"; E = "  -3.04%  " }
    @{ Row = 4; D = "1.007"; E = "  +0.38%  " }
    @{ Row = 5; D = "1.007"; E = "  +0.46%  " }
    @{ Row = 6; D = "307.67"; E = "  -2.41%  " }
    @{ Row = 7; D = "0.4538"; E = $null }
    @{ Row = 8; D = "0.3644"; E = "  -1.90%  " }
    @{ Row = 9; D = "0.07094"; E = "  -3.10%  " }
    @{ Row = 10; D = "0.8701"; E = "  -2.14%  " }
    @{ Row = 12; D = "19.22"; E = "  -4.46%  " }
    @{ Row = 13; D = "1.879.05"; E = "  +0.68%  " }
    @{ Row = 14; D = "5.263"; E = "  -2.54%  " }
    @{ Row = 15; D = "6.324"; E = "  -3.47%  " }
    @{ Row = 16; D = "86.42"; E = "  -5.93%  " }
    @{ Row = 17; D = $null; E = "  +0.68%  " }
    @{ Row = 18; D = "0.000008556"; E = "  -4.73%  " }
    @{ Row = 19; D = $null; E = "  +0.66%  " }
    @{ Row = 20; D = "26.452.57"; E = "  -3.47%  " }
    @{ Row = 21; D = "14.20"; E = "  -4.09%  " }
    @{ Row = 22; D = "4.946"; E = "  -3.63%  " }
    @{ Row = 23; D = "2.069.22"; E = "  +1.53%  " }
    @{ Row = 24; D = $null; E = "  -2.19%  " }
    @{ Row = 25; D = "1.972"; E = "  +1.78%  " }
    @{ Row = 26; D = "150.50"; E = "  -1.08%  " }
    @{ Row = 27; D = "17.85"; E = "  -3.07%  " }
    @{ Row = 28; D = "1.996"; E = "  -2.76%  " }
    @{ Row = 29; D = "113.03"; E = "  -2.76%  " }
    @{ Row = 30; D = "4.858"; E = "  -4.82%  " }
    @{ Row = 31; D = "0.08670"; E = "  -2.06%  " }
    @{ Row = 32; D = "3.111"; E = "  -0.83%  " }
    @{ Row = 33; D = "0.7263"; E = "  -5.80%  " }
    @{ Row = 34; D = "4.420"; E = "  -2.07%  " }
    @{ Row = 35; D = $null; E = "  -5.47%  " }
    @{ Row = 36; D = $null; E = "  +0.83%  " }
    @{ Row = 37; D = "2.495"; E = "  -9.95%  " }
    @{ Row = 38; D = "1.073"; E = "  -0.71%  " }
    @{ Row = 39; D = "0.01905"; E = "  -2.81%  " }
    @{ Row = 40; D = "0.05082"; E = "  -3.14%  " }
    @{ Row = 41; D = "2.848"; E = "  -3.66%  " }
    @{ Row = 42; D = "6.872"; E = "  -2.96%  " }
    @{ Row = 43; D = "0.4890"; E = $null }
    @{ Row = 44; D = "0.1565"; E = "  -4.60%  " }
    @{ Row = 45; D = "8.115"; E = "  -3.61%  " }
    @{ Row = 46; D = "1.008"; E = "  +0.54%  " }
    @{ Row = 47; D = $null; E = "  -4.73%  " }
    @{ Row = 48; D = "101.46"; E = "  -1.71%  " }
    @{ Row = 49; D = "9.912"; E = "  -4.28%  " }
    @{ Row = 50; D = "1.575"; E = "  -4.69%  " }
    @{ Row = 51; D = "0.05988"; E = "  -3.74%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Range("D" + $u.Row)
        $priceIsNumericLike = $u.D -match "^[+-]?\d+(\.\d+)?$"
        if ($priceIsNumericLike) {
            $cell.NumberFormat = "@"
            $cell.Value = $u.D
            $cell.ClearFormats()
        } else {
            $cell.Value = $u.D
        }
    }
    if ($null -ne $u.E) {
        $ws.Range("E" + $u.Row).Value = $u.E
    }
}
